$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7981618046760559
$ws.Range("B1").Value = 1.15644097328186
$ws.Range("C1").Value = 2.259835481643677
$ws.Range("D1").Value = 3.932620763778687
$ws.Range("E1").Value = 1.918848633766174
